$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "65.922.96"
$ws.Range("E2").Value = "  +7.07%  "
$ws.Range("D3").Value = "3.015.97"
$ws.Range("E3").Value = "  +4.26%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "586.38"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +3.25%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "156.39"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +9.31%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("D8").Value = "3.011.27"
$ws.Range("E8").Value = "  +4.16%  "
$ws.Range("E9").Value = "  +2.85%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.01"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +1.39%  "
$ws.Range("E11").Value = "  +7.04%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.453"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +5.47%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000251"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +9.22%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "34.47"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +8.44%  "
$ws.Range("E15").Value = "  +0.64%  "
$ws.Range("D16").Value = "65.896.88"
$ws.Range("E16").Value = "  +7.04%  "
$ws.Range("D17").Value = "3.516.66"
$ws.Range("E17").Value = "  +4.33%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.97"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +6.62%  "
$ws.Range("D19").Value = "3.019.84"
$ws.Range("E19").Value = "  +4.60%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "464.72"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +7.64%  "
$ws.Range("E21").Value = "  +6.35%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.684"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +4.78%  "
$ws.Range("E23").Value = "  +8.12%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "82.35"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +4.00%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "12.53"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +5.55%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.82"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +9.64%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.999"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -0.09%  "
$ws.Range("B29").Value = "ImmutableX"
$ws.Range("C29").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.41"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +18.45%  "
$ws.Range("B30").Value = "NEARProtocol"
$ws.Range("C30").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.91"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +13.23%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0000107"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +1.28%  "
$ws.Range("E32").Value = "  +4.97%  "
$ws.Range("E33").Value = "  +5.39%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "27.04"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +6.03%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.999"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -0.06%  "
$ws.Range("E36").Value = "  +4.46%  "
$ws.Range("B37").Value = "Stacks"
$ws.Range("C37").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.20"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +13.96%  "
$ws.Range("B38").Value = "Filecoin"
$ws.Range("C38").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "5.81"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +8.18%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.04"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +8.47%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "49.24"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +0.71%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "44.90"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +13.17%  "
$ws.Range("E42").Value = "  +8.13%  "
$ws.Range("E43").Value = "  +12.14%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.51"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +3.85%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "392.55"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +14.26%  "
$ws.Range("D46").Value = "2.805.80"
$ws.Range("E46").Value = "  +4.36%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0355"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +6.13%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "134.02"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +0.95%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "23.71"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +10.26%  "
$ws.Range("E51").Value = "  +4.46%  "
